# Add team W/L/T record columns (Wins, Losses, Ties) to the roster sheet.
# New columns: AD = Wins, AE = Losses, AF = Ties.
# Header row (row 1) gets the same bold/bordered/centered style already
# used by the other header cells (e.g. AC1). Data rows 2-55 get the
# team's record as plain numbers: 80 wins, 82 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (from AC1, the last header cell)
# onto the three new header cells so they reuse the same cell style
# instead of minting a new one.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every player row (2 through 55).
$ws.Range("AD2:AD55").Value = 80
$ws.Range("AE2:AE55").Value = 82
$ws.Range("AF2:AF55").Value = 0
